$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add note about error dialogs / showErrorDialogs property ---
$ws.Range("B3").Value = "Note: If no error dialogs are shown in the ARE gui, ensure to set the following property in the fie 'areProperties':`nshowErrorDialogs=1"
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 47.25

# --- Row 7 (ARE_START_3 / start.bat variant) ---
$ws.Range("B7").Value = "Start Are with start.bat (start.sh)"
$ws.Range("C7").Value = "ARE start file: start.bat (start.sh - Linux)"
$ws.Range("D7").Value = "1. Execute ARE start file"

# --- Row 8 (ARE_START_3 / start_debug.bat variant) ---
$ws.Range("B8").Value = "Start Are with start_debug.bat (start_debug.sh)"
$ws.Range("C8").Value = "ARE start file: start_debug.bat (start_debug.sh - Linux)"
$ws.Range("D8").Value = "1. Execute ARE start file"

# --- Row 9 (ARE_START_3 / Are.exe variant) ---
$ws.Range("B9").Value = "Start Are with Are.exe`nTODO: Add Test with ACS connected and ""Download Component Collection"""
$ws.Range("C9").Value = "ARE start file: Are.exe - Windows only"
$ws.Range("D9").Value = "1. Execute ARE start file"

# --- Row 10 (ARE_START_4, damaged model file) ---
$ws.Range("B10").Value = "Use damaged model file autostart_damaged.acs (xml-schema error)"
$ws.Range("C10").Value = "ARE start file: Are.exe (start.sh - Linux)`nmodel file: ARE startup/autostart_damaged.acs"
$ws.Range("D10").Value = "1. Backup original model file models/autostart.xml`n2. Copy model file to models/autostart.xml`n3.Execute ARE start file`n4. Restore original model file"

# --- Row 11 (ARE_START_5, old model file) ---
$ws.Range("B11").Value = "Use old model file Ergo-Kopf-Musik-einfach_v2.5.acs (not up 2 date with bundle_descriptors of used plugins)"
$ws.Range("C11").Value = "ARE start file: Are.exe (start.sh - Linux)`nARE startup/Ergo-Kopf-Musik-einfach_v2.5.acs"
$ws.Range("D11").Value = "1. Backup original model file models/autostart.xml`n2. Copy model file to models/autostart.xml`n3.Execute ARE start file`n4. Restore original model file"
$ws.Range("E11").Value = "ARE should stop starting with an informative error message"

# --- Row 12 (new test case ARE_START_6, Gamepad) ---
$ws.Range("A12").Value = "ARE_START_6"
$ws.Range("B12").Value = "Use old model file GamepadMouse.acs (without plugging in Gamepad)"
$ws.Range("C12").Value = "ARE start file: Are.exe (start.sh - Linux)`nARE startup/GamepadMouse.acs"
$ws.Range("D12").Value = "1. Backup original model file models/autostart.xml`n2. Copy model file to models/autostart.xml`n3. DO NOT plug in Gamepad`n4.Execute ARE start file`n5. Restore original model file"
$ws.Range("E12").Value = "ARE should start with the model running and  show up informative error message"

# --- View / selection changes ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
